# Update annotations for Ying Tang
# - Row 76, column B (politeness_score) becomes a true numeric value (4) instead of text "4".
# - A new row 77 is appended with a new annotation record (politeness_score kept as text "4").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 76: make B76 a genuine number (4) ---
$ws.Cells.Item(76, 2).Value = 4

# --- Row 77: new annotation row ---
$ws.Cells.Item(77, 1).Value = "Ying Tang"

# B77 must stay a text value "4" (not a number), so force the text number
# format before assigning, then restore the default "Normal" style so no
# stray style index is left behind on the cell.
$ws.Cells.Item(77, 2).NumberFormat = "@"
$ws.Cells.Item(77, 2).Value = "4"
$ws.Cells.Item(77, 2).Style = "Normal"

$ws.Cells.Item(77, 3).Value = "I do agree， I have not seen other works， missing，hinder this paper significantly"
$ws.Cells.Item(77, 4).Value = "FBK"
$ws.Cells.Item(77, 5).Value = "OTH"
$ws.Cells.Item(77, 6).Value = "84d24e0a-0b18-4f4c-a441-4ea819712145"
$ws.Cells.Item(77, 7).Value = "rk9kKMZ0-_annotated.xlsx"
$ws.Cells.Item(77, 8).Value = "While the idea is novel and I do agree that I have not seen other works along these lines there are a few things that are missing and hinder this paper significantly."
